$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.533.59"
$ws.Range("E2").Value = "  -4.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.98"
$ws.Range("E3").Value = "  -3.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "273.93"
$ws.Range("E5").Value = "  -8.74%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5010"
$ws.Range("E7").Value = "  -5.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3401"
$ws.Range("E8").Value = "  -8.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.73"
$ws.Range("E9").Value = "  -3.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06627"
$ws.Range("E10").Value = "  -7.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.40"
$ws.Range("E11").Value = "  -9.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7927"
$ws.Range("E12").Value = "  -10.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07823"
$ws.Range("E13").Value = "  -4.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.816.98"
$ws.Range("E14").Value = "  -2.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.006"
$ws.Range("E15").Value = "  -5.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.61"
$ws.Range("E16").Value = "  -6.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.92"
$ws.Range("E18").Value = "  -6.00%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007935"
$ws.Range("E20").Value = "  -6.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.586.25"
$ws.Range("E21").Value = "  -4.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.700"
$ws.Range("E22").Value = "  -5.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.841"
$ws.Range("E23").Value = "  -7.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.082"
$ws.Range("E24").Value = "  -4.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.236"
$ws.Range("E25").Value = "  -1.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.36"
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.656"
$ws.Range("E27").Value = "  -4.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.00"
$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("E29").Value = "  -4.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.236"
$ws.Range("E30").Value = "  -9.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.175"
$ws.Range("E31").Value = "  -9.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08674"
$ws.Range("E32").Value = "  -4.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04758"
$ws.Range("E33").Value = "  -5.08%  "

$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.842"
$ws.Range("E35").Value = "  -3.45%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7104"
$ws.Range("E36").Value = "  -12.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.110"
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.302"
$ws.Range("E39").Value = "  -12.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01821"
$ws.Range("E40").Value = "  -6.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4981"
$ws.Range("E41").Value = "  -18.52%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9307"
$ws.Range("E42").Value = "  -12.31%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.43"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.132"

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.737"
$ws.Range("E46").Value = "  -11.16%  "

$ws.Range("E47").Value = "  -9.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4332"
$ws.Range("E48").Value = "  -17.61%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.177"
$ws.Range("E49").Value = "  -7.68%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.90"
$ws.Range("E50").Value = "  -3.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05899"
$ws.Range("E51").Value = "  -2.64%  "
